$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.034.31"
$ws.Range("E2").Value = "  +1.03%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.092.12"
$ws.Range("E3").Value = "  +0.07%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'540.75"
$ws.Range("E5").Value = "  -1.44%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'137.21"

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.03%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.087.74"
$ws.Range("E8").Value = "  +0.25%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +0.79%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "'6.43"
$ws.Range("E11").Value = "  -2.72%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -0.82%  "

# Row 13 - ShibaInu
$ws.Range("D13").Value = "'0.0000229"
$ws.Range("E13").Value = "  +5.25%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'34.80"
$ws.Range("E14").Value = "  -0.73%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.593.19"
$ws.Range("E15").Value = "  +0.12%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "64.046.56"
$ws.Range("E16").Value = "  +1.00%  "

# Row 17 - TRON
$ws.Range("E17").Value = "  +1.07%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.087.23"
$ws.Range("E18").Value = "  -0.12%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "'6.70"
$ws.Range("E19").Value = "  +0.45%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'485.07"
$ws.Range("E20").Value = "  -0.08%  "

# Row 21 - Chainlink
$ws.Range("D21").Value = "'13.42"
$ws.Range("E21").Value = "  -0.18%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "'0.703"
$ws.Range("E22").Value = "  -0.04%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  -0.90%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'79.75"
$ws.Range("E24").Value = "  +2.49%  "

# Row 25 - InternetComputer(DFINITY)
$ws.Range("D25").Value = "'12.25"
$ws.Range("E25").Value = "  +0.06%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.21%  "

# Row 27 - PancakeSwap
$ws.Range("E27").Value = "  -0.86%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "'8.11"
$ws.Range("E28").Value = "  -1.87%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "'26.42"

# Row 31 - ImmutableX
$ws.Range("E31").Value = "  -2.30%  "

# Row 32 - Mantle
$ws.Range("D32").Value = "'1.14"
$ws.Range("E32").Value = "  +0.94%  "

# Row 33 - OKB
$ws.Range("D33").Value = "'57.93"
$ws.Range("E33").Value = "  -5.32%  "

# Row 34 - Stacks
$ws.Range("E34").Value = "  -6.31%  "

# Row 35 - Bittensor
$ws.Range("D35").Value = "'503.55"

# Row 36 - NEARProtocol
$ws.Range("D36").Value = "'5.37"
$ws.Range("E36").Value = "  +3.80%  "

# Row 37 - Filecoin
$ws.Range("E37").Value = "  +0.91%  "

# Row 38 - Maker
$ws.Range("D38").Value = "3.241.03"
$ws.Range("E38").Value = "  +5.53%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "'0.0400"
$ws.Range("E39").Value = "  -0.56%  "

# Row 40 - Hedera
$ws.Range("D40").Value = "'0.0798"
$ws.Range("E40").Value = "  +0.49%  "

# Row 41 - Kaspa
$ws.Range("E41").Value = "  +0.68%  "

# Rows 42/43 - Cosmos and dogwifhat swap ranking order
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.68"
$ws.Range("E42").Value = "  +0.25%  "

$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "'8.14"
$ws.Range("E43").Value = "  +0.29%  "

# Row 44 - TheGraph
$ws.Range("D44").Value = "'0.256"
$ws.Range("E44").Value = "  +0.51%  "

# Row 45 - USDe
$ws.Range("E45").Value = "  +0.07%  "

# Rows 46/47 - Fetch.AI and Monero swap ranking order
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.05"
$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'122.76"
$ws.Range("E47").Value = "  +0.95%  "

# Row 48 - PEPE
$ws.Range("E48").Value = "  +5.36%  "

# Row 49 - InjectiveProtocol
$ws.Range("D49").Value = "'24.71"
$ws.Range("E49").Value = "  +1.45%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  +1.71%  "

# Row 51 - CoreDAO
$ws.Range("D51").Value = "'2.41"
$ws.Range("E51").Value = "  +3.14%  "
